# Update "Förändrad" (changed) date column C for rows 2-103 from 45188 (2023-09-19)
# to 45189 (2023-09-20), i.e. advance each value by one day.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$lastRow = 103
$firstRow = 2

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)  # Column C
    if ($cell.Value2 -eq 45188) {
        $cell.Value2 = 45189
    }
}
